# Generate Report for Handoff
# Rotate the localization-status report from the old source GUID
# (b9270cb6-8e94-4286-be2b-bb852b889057) to the new one
# (bd18bc81-4910-4935-8288-b26f75fb1fb1), refreshing the associated
# handoff/handback hashes and timestamps on every sheet.

$wb = $excel.ActiveWorkbook

$oldGuid = "b9270cb6-8e94-4286-be2b-bb852b889057"
$newGuid = "bd18bc81-4910-4935-8288-b26f75fb1fb1"

# The hyperlinks on every sheet still resolve to the original file at the
# original commit - only the cell text / display label is rotated.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/26fbf973aebf0323d41831a8743e2e911930b232/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# "Overview" sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newGuid.md"

$newDisplayB2 = "e2e\$newGuid.md"
$ws.Range("B2").Value = $newDisplayB2
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $linkAddress, "", "", $newDisplayB2)

$ws.Range("G2").Value = "2016-09-03 23:01:47"

# ---------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$newDisplayA2 = "$newGuid.md"
$ws.Range("A2").Value = $newDisplayA2
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $linkAddress, "", "", $newDisplayA2)

$ws.Range("G2").Value = "$newGuid.b3df7d475817a3cc48069f9098289acc38d3ab2f.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-03 23:01:42"

# ---------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$newDisplayA2 = "$newGuid.md"
$ws.Range("A2").Value = $newDisplayA2
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $linkAddress, "", "", $newDisplayA2)

$ws.Range("G2").Value = "$newGuid.b3df7d475817a3cc48069f9098289acc38d3ab2f.de-de.xlf"
$ws.Range("H2").Value = "2016-09-03 23:01:47"
